$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 22, pushing the "Next steps" / "Merge to NCS"
# block (old rows 22-24) down to rows 24-26.
$ws.Rows("22:23").Insert()

# New "Build command" bold/yellow header cell at A18 (row 18 already holds the
# "Replace in Wi-Fi and BLE coex sample." text in column C).
$ws.Range("A18").Value = "Build command"
$ws.Range("A18").Font.Bold = $true
$ws.Range("A18").Interior.Color = 65535

# New "TODO" label in A22 (now an empty row after the insert above).
$ws.Range("A22").Value = "TODO"

# Page setup: paper size / orientation were touched in the saved file.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the view state that the diff shows (selection can be reproduced;
# scroll position is best-effort).
$ws.Range("E14").Select()
